# Update the Controller worksheet:
#  - D2 becomes the string "sachin" (was the number 1)
#  - Add rows 3-7 with new test data (Registration.xlsx / RestTemplate.xlsx rows)
#  - Column B gets a wider custom width
#  - Selection moves to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2 switches from a literal number to the text "sachin"
$ws.Range("D2").Value = "sachin"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Registration.xlsx"
$ws.Range("C3").Value = "REG"
$ws.Range("D3").Value = 2
$ws.Range("G3").Value = "chrome"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Registration.xlsx"
$ws.Range("C4").Value = "REG"
$ws.Range("D4").Value = 3
$ws.Range("G4").Value = "chrome"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "RestTemplate.xlsx"
$ws.Range("C5").Value = "REG"
$ws.Range("D5").Value = 4
$ws.Range("F5").Value = 1

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "RestTemplate.xlsx"
$ws.Range("C6").Value = "REG"
$ws.Range("D6").Value = 5
$ws.Range("F6").Value = 1

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "RestTemplate.xlsx"
$ws.Range("C7").Value = "REG"
$ws.Range("D7").Value = 6
$ws.Range("F7").Value = 1

# Widen column B (32 characters wide, enough to fit "Registration.xlsx" / "RestTemplate.xlsx")
# and move the active selection to D2, matching the saved view state
$ws.Columns.Item(2).ColumnWidth = 32
$ws.Range("D2").Select()
